$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - A2, B2, C2 change; D2, E2, F2 unchanged
$ws.Range("A2").Value = "trainingimages/18_popata"
$ws.Range("B2").Value = "pngimages/18_donut.png"
$ws.Range("C2").Value = "trainingimages/05_titopo"

# Row 3 - all cells change
$ws.Range("A3").Value = "trainingimages/23_patoko"
$ws.Range("B3").Value = "pngimages/23_lemon.png"
$ws.Range("C3").Value = "trainingimages/06_titoka"
$ws.Range("D3").Value = "pngimages/06_tent.png"
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = -0.5

# Row 4 - all cells change
$ws.Range("A4").Value = "trainingimages/04_kitoti"
$ws.Range("B4").Value = "pngimages/04_ladder.png"
$ws.Range("C4").Value = "trainingimages/21_papika"
$ws.Range("D4").Value = "pngimages/21_cheese.png"
$ws.Range("E4").Value = -0.5
$ws.Range("F4").Value = 0.5
